# "class status update working"
# Adds a new field "studentclassstatus" (varchar(40)) to the code-generator
# workbook:
#   - Sheet1: new row 40, mirroring the code-gen formulas of row 39, for the
#     new field name held in A40.
#   - Sheet3 (the field-type lookup sheet, stored as worksheets/sheet2.xml):
#     new row 18 describing the new column's SQL type (varchar(40)).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------
# Sheet3 (lookup table of field name -> sql type / char columns) - row 18
# ---------------------------------------------------------------------
$ws3.Range("B18").Value = "studentclassstatus"
$ws3.Range("C18").Value = "varchar(40)"
$ws3.Range("D18").Value = "CHARACTER"
$ws3.Range("E18").Value = "SET"
$ws3.Range("F18").Value = "utf8"
$ws3.Range("G18").Value = "DEFAULT"
$ws3.Range("H18").Value = "NULL,"
$ws3.Range("I18").Value = "s"

# ---------------------------------------------------------------------
# Sheet1 - row 40, the new "studentclassstatus" field, following the
# same per-column formula pattern as the preceding rows (35-39).
# ---------------------------------------------------------------------
$ws1.Range("A40").Value = "studentclassstatus"

$ws1.Range("B40").Formula = '="$sql .= """&" t."&TRIM(A40)&" = :"&TRIM(A40)&", "&""";"'
$ws1.Range("C40").Formula = '="$studentclass->"&TRIM(A40)&","'
$ws1.Range("D40").Formula = '="$"&TRIM(A40)&" = $app->request->put(''"&TRIM(A40)&"'');"'
$ws1.Range("E40").Formula = '=VLOOKUP(A40,Sheet3!B:I,8,FALSE)'
$ws1.Range("F40").Formula = '=F39&E40'
$ws1.Range("G40").Formula = '= "error_log( print_R($"&TRIM(A40)&", TRUE ));"'
$ws1.Range("H40").Formula = '="$response["""&TRIM(A40)&"""] = $result["""&TRIM(A40)&"""];"'
$ws1.Range("I40").Formula = '="$"&TRIM(A40)&" = $studentclass->"&TRIM(A40)&";"'
$ws1.Range("J40").Formula = '="$"&TRIM(A40)&","'
$ws1.Range("K40").Formula = '="$res[""sc_"&TRIM(A40)&"""] = $sc_"&TRIM(A40)&";"'
$ws1.Range("L40").Formula = '="t."&TRIM(A40)&","'
$ws1.Range("M40").Formula = '="$sc_"&TRIM(A40)&","'

# ---------------------------------------------------------------------
# Selection / active-sheet state, matching the author's final view:
# Sheet3 (B:I lookup) scrolled to its new row, then Sheet1 left active
# with the new row's E:F cells selected.
# ---------------------------------------------------------------------
$ws3.Activate() | Out-Null
$ws3.Range("I18").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("E40:F40").Select() | Out-Null
